$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title (A1): "October 2016 and 2015" -> "November 2016 and 2015" ---
$ws.Range("A1").Value = "Table 4.13.B. Average Cost of Natural Gas Delivered for Electricity Generation by State, (Year-to-Date) November 2016 and 2015"

# --- Column group headers (row 4): "October 20xx YTD" -> "November 20xx YTD" ---
$ws.Range("B4").Value = "November 2016 YTD"
$ws.Range("C4").Value = "November 2015 YTD"
$ws.Range("E4").Value = "November 2016 YTD"
$ws.Range("F4").Value = "November 2015 YTD"
$ws.Range("G4").Value = "November 2016 YTD"
$ws.Range("H4").Value = "November 2015 YTD"

# --- Data cell updates (rows 5-66), state-by-state refreshed monthly figures ---
$updates = @(
    @("B5", 3.08),
    @("C5", 4.34),
    @("D5", -0.29),
    @("F5", 3.93),
    @("G5", 3.08),
    @("H5", 4.35),
    @("B6", 3.39),
    @("C6", 4.57),
    @("D6", -0.26),
    @("G6", 3.39),
    @("H6", 4.57),
    @("B8", 2.89),
    @("C8", 4.3),
    @("D8", -0.33),
    @("G8", 2.88),
    @("H8", 4.31),
    @("F9", 4.91),
    @("C10", 3.69),
    @("H10", 3.69),
    @("B12", 2.08),
    @("C12", 2.99),
    @("D12", -0.3),
    @("E12", 2.5),
    @("F12", 3.71),
    @("G12", 2.03),
    @("H12", 2.9),
    @("B13", 2.01),
    @("C13", 2.97),
    @("D13", -0.32),
    @("G13", 2.01),
    @("H13", 2.97),
    @("B14", 2.49),
    @("C14", 3.45),
    @("D14", -0.28),
    @("E14", 2.5),
    @("F14", 3.71),
    @("G14", 2.49),
    @("H14", 3.36),
    @("B15", 1.74),
    @("C15", 2.5),
    @("D15", -0.3),
    @("G15", 1.74),
    @("H15", 2.5),
    @("B16", 2.57),
    @("C16", 2.89),
    @("D16", -0.11),
    @("F16", 3.02),
    @("H16", 2.78),
    @("C17", "W"),
    @("D17", "W"),
    @("F17", 3.81),
    @("H17", "W"),
    @("E18", 2.83),
    @("F18", 3.01),
    @("C19", 3.23),
    @("D19", -0.18),
    @("F19", 3.22),
    @("G19", 2.54),
    @("H19", 3.24),
    @("B20", 2.15),
    @("C20", 2.33),
    @("D20", -0.077),
    @("F20", 2.44),
    @("G20", 2.14),
    @("H20", 2.29),
    @("E21", 2.72),
    @("F21", 3.24),
    @("B22", 2.86),
    @("E22", 2.86),
    @("F22", 3.51),
    @("G22", 2.84),
    @("B23", 2.58),
    @("C23", 3.12),
    @("E23", 2.58),
    @("F23", 3.12),
    @("B24", 3.27),
    @("E24", 3.27),
    @("E25", 2.99),
    @("F25", 3.73),
    @("E26", 2.77),
    @("F26", 3.32),
    @("B27", 3.05),
    @("C27", 3.72),
    @("E27", 3.05),
    @("F27", 3.72),
    @("C28", 8.42),
    @("F28", 8.42),
    @("B29", 2.42),
    @("C29", 3.24),
    @("D29", -0.25),
    @("E29", 2.42),
    @("F29", 3.24),
    @("B30", 3.39),
    @("C30", 4.02),
    @("D30", -0.16),
    @("E30", 3.49),
    @("F30", 4.16),
    @("G30", 2.58),
    @("H30", 2.91),
    @("B33", 3.73),
    @("C33", 4.35),
    @("D33", -0.14),
    @("E33", 3.75),
    @("F33", 4.37),
    @("G33", 2.85),
    @("H33", 2.84),
    @("B34", 2.9),
    @("C34", 3.23),
    @("D34", -0.1),
    @("E34", 2.96),
    @("F34", 3.27),
    @("G34", 2.69),
    @("H34", 3.08),
    @("C35", 3.96),
    @("D35", -0.29),
    @("H35", 3.96),
    @("F36", 4.72),
    @("F37", 3.43),
    @("C38", 3.45),
    @("E38", 2.97),
    @("F38", 3.91),
    @("H38", 2.18),
    @("F39", 2.82),
    @("B40", 2.74),
    @("C40", 2.98),
    @("D40", -0.081),
    @("E40", 2.75),
    @("F40", 2.96),
    @("G40", 2.73),
    @("H40", 3.01),
    @("E41", 2.86),
    @("F41", 3.07),
    @("F42", 3.53),
    @("E43", 2.72),
    @("F43", 2.9),
    @("B44", 2.51),
    @("C44", 2.76),
    @("D44", -0.091),
    @("E44", 2.51),
    @("F44", 2.76),
    @("C45", 2.88),
    @("D45", -0.11),
    @("E45", 2.68),
    @("F45", 2.99),
    @("H45", 2.8),
    @("F46", 3.35),
    @("B47", 2.6),
    @("C47", "W"),
    @("E47", 2.65),
    @("F47", 2.98),
    @("G47", 2.37),
    @("H47", "W"),
    @("E48", 2.68),
    @("F48", 3.05),
    @("B49", 2.55),
    @("C49", 2.84),
    @("D49", -0.1),
    @("E49", 2.65),
    @("F49", 2.95),
    @("H49", 2.81),
    @("B50", 2.91),
    @("E50", 2.92),
    @("F50", 3.28),
    @("G50", 2.76),
    @("C51", 3.34),
    @("E51", 3.11),
    @("F51", 3.46),
    @("H51", 2.99),
    @("C52", 3.57),
    @("E52", 3.06),
    @("F52", 3.44),
    @("H52", 4.31),
    @("B53", 2.84),
    @("C53", 2.94),
    @("D53", -0.034),
    @("E53", 2.84),
    @("F53", 2.94),
    @("F54", 2.31),
    @("B55", 2.85),
    @("C55", 3.23),
    @("D55", -0.12),
    @("E55", 2.85),
    @("F55", 3.23),
    @("B56", 2.84),
    @("C56", 3.15),
    @("D56", -0.098),
    @("E56", 2.84),
    @("F56", 3.15),
    @("E57", 2.58),
    @("F57", 2.96),
    @("B58", 7.95),
    @("E58", 7.95),
    @("F58", 4.64),
    @("B59", 2.93),
    @("C59", 3.24),
    @("D59", -0.096),
    @("E59", 3.2),
    @("F59", 3.5),
    @("G59", 2.69),
    @("H59", 3),
    @("B60", 3.05),
    @("C60", 3.33),
    @("D60", -0.084),
    @("E60", 3.44),
    @("F60", 3.67),
    @("G60", 2.78),
    @("H60", 3.07),
    @("E61", 2.29),
    @("F61", 2.78),
    @("E62", 3.33),
    @("F62", 3.49),
    @("B63", 6.57),
    @("C63", 5.37),
    @("E63", 6.57),
    @("F63", 5.37),
    @("B64", 6.57),
    @("C64", 5.37),
    @("E64", 6.57),
    @("F64", 5.37),
    @("B66", 2.81),
    @("C66", 3.32),
    @("D66", -0.15),
    @("E66", 3.08),
    @("F66", 3.57),
    @("G66", 2.45),
    @("H66", 3)
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}

Write-Host "Applied $($updates.Count) data cell updates plus title/header changes."
